# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (aggregate/industry row) ---
$ws.Range("D2").Value = 0.423
$ws.Range("K2").Value = 1.26
$ws.Range("L2").Value = 0.1305699481865285

$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0

$ws.Range("U2").Value = 39.2
$ws.Range("V2").Value = 0.9053117782909932
$ws.Range("W2").Value = 0.02709677419354839
$ws.Range("X2").Value = 0.03130610771243929
$ws.Range("Y2").Value = -0.004209333518890902
$ws.Range("Z2").Value = -1.66092943201377
$ws.Range("AA2").Value = -0

$ws.Range("AB2").Value = 0.03080994059358636
$ws.Range("AC2").Value = -0.03080994059358636
$ws.Range("AD2").Value = 1.4
$ws.Range("AF2").Value = 1.4
$ws.Range("AG2").Value = -37.8
$ws.Range("AH2").Value = 0.03131991051454139
$ws.Range("AI2").Value = 0.02845528455284553
$ws.Range("AJ2").Value = -6.872727272727283
$ws.Range("AK2").Value = -3.780000000000003

# --- Row 3 (individual company row) ---
$ws.Range("B3").Value = "BH Mubasher Financial Services P.S.C (DFM:BHMUBASHER)"

$ws.Range("D3").Value = 0.423
$ws.Range("K3").Value = 1.26
$ws.Range("L3").Value = 0.1305699481865285

$ws.Range("O3").Value = -0
$ws.Range("R3").Value = -0

$ws.Range("U3").Value = 39.2
$ws.Range("V3").Value = 0.9053117782909932
$ws.Range("W3").Value = 0.02709677419354839
$ws.Range("X3").Value = 0.03130610771243929
$ws.Range("Y3").Value = -0.004209333518890902
$ws.Range("Z3").Value = -1.66092943201377
$ws.Range("AA3").Value = -0

$ws.Range("AB3").Value = 0.03080994059358636
$ws.Range("AC3").Value = -0.03080994059358636
$ws.Range("AD3").Value = 1.4
$ws.Range("AF3").Value = 1.4
$ws.Range("AG3").Value = -37.8
$ws.Range("AH3").Value = 0.03131991051454139
$ws.Range("AI3").Value = 0.02845528455284553
$ws.Range("AJ3").Value = -6.872727272727283
$ws.Range("AK3").Value = -3.780000000000003
